$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text (preserves values like "1.220", "14.20", "0.3380"
# that would otherwise be auto-converted to numbers by Excel, losing formatting).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.461.82'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.723.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.97'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3733'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.14'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3380'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.41%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07430'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.388'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.99'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.025'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +6.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.721.46'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.55%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06661'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.92'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.21%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.51'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.168'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.87%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.69'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '26.449.78'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.73%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.458'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.407'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +20.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.381'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.65%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '150.69'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.39'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.913.68'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '131.21'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.095'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.951'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08615'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.687'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.17%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.350'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02331'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.95%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06201'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.91%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2151'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.378'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.220'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6197'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.20'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.60%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.891'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.05%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6007'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.42'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.039'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07165'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '76.81'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.52%  '
